$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "SAT Jan 20" + " 10:34:36 PST 2018" (two runs) were merged
# into a single run "SAT Jan 20 10:34:36 PST 2018" by a Find & Replace.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("SAT Jan 20 10:34:36 PST 2018", $false, $false, $false, $false, $false, $true, 1, $false, "SAT Jan 20 10:34:36 PST 2018", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: append a new purchase-record block (28/01/2018 MAHADEVA
# CHICK IN) right after the last existing record in the document.
# ---------------------------------------------------------------------

# Locate the end of the very last "- CASH" occurrence in the document.
$rng = $d.Content
$lastEnd = -1
while ($rng.Find.Execute("- CASH", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastEnd = $rng.End
    $rng.Start = $rng.End
    $rng.End = $d.Content.End
}

# The paragraph right after "- CASH" is an existing blank paragraph;
# the new block is inserted right after that one (i.e. right after its
# paragraph mark).
$insertPos = $lastEnd + 1
$ins = $d.Range($insertPos, $insertPos)

# Blank paragraph
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "SAT Jan 27 10:59:12 PST 2018"
$ins.InsertAfter("SAT Jan 27")
$ins.Collapse(0)
$ins.InsertAfter(" 10:59:12 PST 2018")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Person Name" ... - MAHADEVA
$ins.InsertAfter("Person Name")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- MAHADEVA")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Bill number" ... - 2790
$ins.InsertAfter("Bill number")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 2790")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# Separator line
$ins.InsertAfter("---------------------------------------------------------------")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Item Name" ... - CARROT
$ins.InsertAfter("Item Name")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- CARROT")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Number of Pockets" ... - 2
$ins.InsertAfter("Number of Pockets")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 2")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Number of KGs" ... - 178
$ins.InsertAfter("Number of KGs")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 178")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Rate" ... - 10
$ins.InsertAfter("Rate")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 10")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Total Price" ... - 1780.0
$ins.InsertAfter("Total Price")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 1780.0")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Amount Received" ... - 3471  (red paragraph)
$amountReceivedStart = $ins.Start
$ins.InsertAfter("Amount Received")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 3471")
$ins.Collapse(0)
$amountReceivedEnd = $ins.Start
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Amount balance" ... - 1780.0  (bold paragraph)
$amountBalanceStart = $ins.Start
$ins.InsertAfter("Amount balance")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- 1780.0")
$ins.Collapse(0)
$amountBalanceEnd = $ins.Start
$ins.InsertAfter("`r")
$ins.Collapse(0)

# "Amount Received mode" ... - CASH
$ins.InsertAfter("Amount Received mode")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("`t")
$ins.Collapse(0)
$ins.InsertAfter("- CASH")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)

# Trailing blank paragraph
$ins.InsertAfter("`r")
$ins.Collapse(0)

# Apply character formatting to the two special paragraphs.
$redRange = $d.Range($amountReceivedStart, $amountReceivedEnd)
$redRange.Font.Color = 255

$boldRange = $d.Range($amountBalanceStart, $amountBalanceEnd)
$boldRange.Font.Bold = 1
